$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burndown Chart")

# Fix the ideal burndown formula: denominator changes from 15 to 3
$ws.Range("E13").Formula = "=`$D`$13-(`$D`$13/3*1)"
$ws.Range("F13").Formula = "=`$D`$13-(`$D`$13/3*2)"
$ws.Range("G13").Formula = "=`$D`$13-(`$D`$13/3*3)"

# Update the selected cell to match the new state left after the edit
$ws.Range("L21").Select()
